# Revert "Worked on BOMs" — restore the previous SparkFun-based crimp/housing/
# header parts (rows 4 & 5), clear out the transistor line (row 7), fix up the
# derived totals, hyperlinks, column width and active-cell selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 4 (Header / Molex 5 pin polar crimp housing -> SparkFun part) -----
$ws.Range("D4").Style = "Normal"
$ws.Range("D4").Value = "Molex"
$ws.Range("E4").Style = "Normal"
$ws.Range("E4").Value = "22-23-2051"
$ws.Range("F4").Value = "Sparkfun"
$ws.Range("G4").Value = "PRT-08230"
$ws.Range("H4").Value = "https://www.sparkfun.com/products/8230"
$ws.Range("I4").Value = 0.45

# --- Row 5 (Housing / Molex 5 pin polar crimp housing -> SparkFun part) ----
$ws.Range("D5").Style = "Normal"
$ws.Range("D5").Value = "Molex"
$ws.Range("E5").ClearContents()
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Value = "Sparkfun"
$ws.Range("G5").Value = "PRT-08098"
$ws.Range("H5").Value = "https://www.sparkfun.com/products/8098"
$ws.Range("I5").Value = 0.45

# --- Row 7 (drop the Transistor line entirely) ------------------------------
$ws.Range("A7:J7").ClearContents()
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"

# --- Hyperlinks: drop the digikey/transistor links, restore the SparkFun ---
# --- links for the revived rows 4/5, keep the untouched ones for 2/3/6 -----
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("G2"), "http://www.adafruit.com/products/572")
$ws.Hyperlinks.Add($ws.Range("G3"), "https://www.sparkfun.com/products/8619?")
$ws.Hyperlinks.Add($ws.Range("G4"), "https://www.sparkfun.com/products/8230")
$ws.Hyperlinks.Add($ws.Range("H4"), "https://www.sparkfun.com/products/8230")
$ws.Hyperlinks.Add($ws.Range("H5"), "https://www.sparkfun.com/products/8098")
$ws.Hyperlinks.Add($ws.Range("G5"), "https://www.sparkfun.com/products/8098")
$ws.Hyperlinks.Add($ws.Range("H6"), "https://www.sparkfun.com/products/8100?")
$ws.Hyperlinks.Add($ws.Range("G6"), "https://www.sparkfun.com/products/8100?")

# --- Column H got a bit wider once the long digikey URLs were gone ---------
$ws.Columns("H").ColumnWidth = 81.5

# --- Window got taller/shorter in the original author's session -----------
$win = $excel.ActiveWindow
$win.Top = 6240
$win.Height = 6300
$win.Left = -15
$win.Width = 25230

# --- Selection moved off the old A5:I5 block onto J11 ----------------------
$ws.Activate()
$ws.Range("J11").Select()
